$d = $word.ActiveDocument

function Locate($text) {
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r
}

# =====================================================================
# PART A: Thesis paragraph -- remove the old "_GoBack" bookmark that
# wraps "Can the problems faced by the Boeing 787 ... innovative " and
# let the two runs that used to be split only by the bookmark markers
# coalesce (this mirrors Word's own run clean-up once the bookmark
# that forced the split is gone).
# =====================================================================

$rAnchor = Locate("hesis: Can the problems faced by the Boeing 787")
$posT       = $rAnchor.Start        # boundary  T | "hesis: "
$posCan     = $posT + 7             # boundary  "hesis: " | "Can the problems..."
$posBoeing  = $posCan + 29          # boundary  "...the" | " Boeing"
$pos787     = $posBoeing + 7        # boundary  " Boeing" | " 787..."

# Keep the surrounding boundaries intact while we coalesce across the
# two points that the old bookmark used to occupy.
$d.Bookmarks.Add("ZZ_BarrierT", $d.Range($posT, $posT)) | Out-Null
$d.Bookmarks.Add("ZZ_BarrierBoeingL", $d.Range($posBoeing, $posBoeing)) | Out-Null
$d.Bookmarks.Add("ZZ_BarrierBoeingR", $d.Range($pos787, $pos787)) | Out-Null

$d.Bookmarks("_GoBack").Delete()

# Nudge the run-list to rebuild (identical text restored right away)
# around the first former-bookmark point ...
$midA = $posCan + 3
$rA = $d.Range($midA, $midA)
$rA.InsertAfter("z")
$d.Range($midA, $midA + 1).Text = ""

# ... and around the second former-bookmark point.
$rSupply = Locate("supply chain strategy?")
$midB = $pos787 + 3
$rB = $d.Range($midB, $midB)
$rB.InsertAfter("z")
$d.Range($midB, $midB + 1).Text = ""

$d.Bookmarks("ZZ_BarrierT").Delete()
$d.Bookmarks("ZZ_BarrierBoeingL").Delete()
$d.Bookmarks("ZZ_BarrierBoeingR").Delete()

# =====================================================================
# PART B: Personal website URL -- insert "-mohamed-" in the middle of
# "https://amineaboussalah.github.io/", producing
# "https://amine-mohamed-aboussalah.github.io/" and leaving the
# "_GoBack" bookmark at the new edit location (right after the typed
# text), matching Word's own behaviour of tracking the last edit spot.
# =====================================================================

$rUrl = Locate("https://amine")
$rUrl.Collapse(0)
$d.Bookmarks.Add("ZZ_BarrierAmine", $rUrl) | Out-Null

$rUrl.InsertAfter("-mohamed-")

$rGoBack = Locate("-mohamed-")
$rGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rGoBack) | Out-Null

$d.Bookmarks("ZZ_BarrierAmine").Delete()
